# Adds a "LineItems" / "LineItems - Formatted" pair of sheets (mirroring the
# existing Items / Items - Formatted pair) and reshapes the "Simple Fields"
# sheets' header/data row to the new invoice-summary layout.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rewrite the "Simple Fields" and "Simple Fields - Formatted" sheets.
# ---------------------------------------------------------------------------

$simple = $wb.Worksheets.Item("Simple Fields")
$simpleFormatted = $wb.Worksheets.Item("Simple Fields - Formatted")

# Shared header row for both sheets.
$simple.Range("A1").Value = "Date"
$simple.Range("B1").Value = "Due Date"
$simple.Range("C1").Value = "Invoice#"
$simple.Range("D1").Value = "Customer ID"
$simple.Range("E1").Value = "Supplier"
$simple.Range("F1").Value = "Subtotal"
$simple.Range("G1").Value = "Tax Rate"
$simple.Range("H1").Value = "Tax Due"
$simple.Range("I1").Value = "Total"
$simple.Range("J1").Value = "LineItems"
$simple.Range("K1").Value = "Items"

$simpleFormatted.Range("A1").Value = "Date"
$simpleFormatted.Range("B1").Value = "Due Date"
$simpleFormatted.Range("C1").Value = "Invoice#"
$simpleFormatted.Range("D1").Value = "Customer ID"
$simpleFormatted.Range("E1").Value = "Supplier"
$simpleFormatted.Range("F1").Value = "Subtotal"
$simpleFormatted.Range("G1").Value = "Tax Rate"
$simpleFormatted.Range("H1").Value = "Tax Due"
$simpleFormatted.Range("I1").Value = "Total"
$simpleFormatted.Range("J1").Value = "LineItems"
$simpleFormatted.Range("K1").Value = "Items"

# Data row - "Simple Fields" keeps the short-form date, "- Formatted" uses ISO.
$simple.Range("A2").Value = "4/7/2019"
$simple.Range("B2").Value = "5/22/2019"
$simple.Range("C2").Value = "850888"
$simple.Range("D2").Value = "A700"
$simple.Range("E2").Value = "Tiefland Glass AG"
$simple.Range("F2").Value = "5,000.00"
$simple.Range("G2").Value = ""
$simple.Range("H2").Value = "500.00"
$simple.Range("I2").Value = "5,500.00"
$simple.Range("J2").Value = "table"
$simple.Range("K2").Value = "table"

$simpleFormatted.Range("A2").Value = "2019-04-07"
$simpleFormatted.Range("B2").Value = "2019-05-22"
$simpleFormatted.Range("C2").Value = "850888"
$simpleFormatted.Range("D2").Value = "A700"
$simpleFormatted.Range("E2").Value = "Tiefland Glass AG"
$simpleFormatted.Range("F2").Value = "5,000.00"
$simpleFormatted.Range("G2").Value = ""
$simpleFormatted.Range("H2").Value = "500.00"
$simpleFormatted.Range("I2").Value = "5,500.00"
$simpleFormatted.Range("J2").Value = "table"
$simpleFormatted.Range("K2").Value = "table"

# ---------------------------------------------------------------------------
# 2. Add the new "LineItems" / "LineItems - Formatted" sheets at the end,
#    mirroring the Items / Items - Formatted sheets with different columns.
# ---------------------------------------------------------------------------

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lineItems = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$lineItems.Name = "LineItems"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lineItemsFormatted = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$lineItemsFormatted.Name = "LineItems - Formatted"

$lineItems.Range("A1").Value = "Description"
$lineItems.Range("B1").Value = "Unit Price"
$lineItems.Range("C1").Value = "QTY"
$lineItems.Range("D1").Value = "Taxed"
$lineItems.Range("E1").Value = "Amount"

$lineItems.Range("A2").Value = "Professional services"
$lineItems.Range("B2").Value = "5,000.00"
$lineItems.Range("C2").Value = "1"
$lineItems.Range("D2").Value = "X"
$lineItems.Range("E2").Value = "5,000.00"

$lineItemsFormatted.Range("A1").Value = "Description"
$lineItemsFormatted.Range("B1").Value = "Unit Price"
$lineItemsFormatted.Range("C1").Value = "QTY"
$lineItemsFormatted.Range("D1").Value = "Taxed"
$lineItemsFormatted.Range("E1").Value = "Amount"

$lineItemsFormatted.Range("A2").Value = "Professional services"
$lineItemsFormatted.Range("B2").Value = "5,000.00"
$lineItemsFormatted.Range("C2").Value = "1"
$lineItemsFormatted.Range("D2").Value = "X"
$lineItemsFormatted.Range("E2").Value = "5,000.00"
